$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean slate: remove old header/data content & formatting ---
$ws.Cells.Clear()

# --- Header text values (row 1) ---
$ws.Range("A1").Value = "Code"
$ws.Range("B1").Value = "Biometric Id"
$ws.Range("C1").Value = "First Name"
$ws.Range("D1").Value = "Last Name"
$ws.Range("E1").Value = "Father/Husband Name"
$ws.Range("F1").Value = "Bloodgroup"
$ws.Range("G1").Value = "Email"
$ws.Range("H1").Value = "Mobile"
$ws.Range("I1").Value = "User Type"
$ws.Range("J1").Value = "Employment Type"
$ws.Range("K1").Value = "Contractor Code"
$ws.Range("L1").Value = "Contractor Name"
$ws.Range("M1").Value = "Membership Date"
$ws.Range("N1").Value = "Supervisor Code"
$ws.Range("O1").Value = "Designation Code"
$ws.Range("P1").Value = "Designation Name"
$ws.Range("Q1").Value = "Department Code"
$ws.Range("R1").Value = "Department Name"
$ws.Range("S1").Value = "Joining Date"
$ws.Range("T1").Value = "Gender"
$ws.Range("U1").Value = "Birthday"
$ws.Range("V1").Value = "Address Line 1"
$ws.Range("W1").Value = "Address Line 2"
$ws.Range("X1").Value = "City"
$ws.Range("Y1").Value = "District"
$ws.Range("Z1").Value = "State"
$ws.Range("AA1").Value = "Pincode"
$ws.Range("AB1").Value = "PF NO"
$ws.Range("AC1").Value = "ESI NO"
$ws.Range("AD1").Value = "UAN NO"
$ws.Range("AE1").Value = "PAN"
$ws.Range("AF1").Value = "Account No"
$ws.Range("AG1").Value = "Account Holder Name"
$ws.Range("AH1").Value = "IFSC"
$ws.Range("AI1").Value = "Bank"
$ws.Range("AJ1").Value = "Branch"
$ws.Range("AK1").Value = "Aadhaar"
$ws.Range("AL1").Value = "Nominee Name"
$ws.Range("AM1").Value = "Nominee Relation"

# --- Formatting ---
# Cells with border + center/center + wrap text (most headers)
$rngWrap = $ws.Range("A1:AA1")
$rngWrap.Borders.LineStyle = 1
$rngWrap.HorizontalAlignment = -4108
$rngWrap.VerticalAlignment = -4108
$rngWrap.WrapText = $true

# Cells with border + center/center, no wrap (PF NO..Nominee Relation block)
$rngNoWrap = $ws.Range("AB1:AM1")
$rngNoWrap.Borders.LineStyle = 1
$rngNoWrap.HorizontalAlignment = -4108
$rngNoWrap.VerticalAlignment = -4108
$rngNoWrap.WrapText = $false

# UAN NO / Account No headers additionally carry a "0" integer number format
$ws.Range("AD1").NumberFormat = "0"
$ws.Range("AF1").NumberFormat = "0"

# Trailing blank header cell (AN1): centered, no border, no wrap
$rngBlank = $ws.Range("AN1")
$rngBlank.HorizontalAlignment = -4108
$rngBlank.VerticalAlignment = -4108
$rngBlank.WrapText = $false

# --- Column widths ---
$ws.Range("A1:B1").EntireColumn.ColumnWidth = 9.0
$ws.Range("C1:G1").EntireColumn.ColumnWidth = 21.6665
$ws.Range("H1:L1").EntireColumn.ColumnWidth = 18.333
$ws.Columns("M").ColumnWidth = 10.0
$ws.Range("AB1:AC1").EntireColumn.ColumnWidth = 13.833
$ws.Columns("AD").ColumnWidth = 12.0
$ws.Columns("AE").ColumnWidth = 20.833
$ws.Columns("AF").ColumnWidth = 20.833
$ws.Columns("AG").ColumnWidth = 17.833
$ws.Columns("AM").ColumnWidth = 15.1665

# --- Sheet view / selection state ---
$ws.Range("AH4").Select()
$excel.ActiveWindow.ScrollColumn = 30
$excel.ActiveWindow.ScrollRow = 1
